# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E22) listed the periods in descending
# order (2110, 2109, ..., 2104). The update re-orders them in ascending
# order (2104, 2105, ..., 2110), and keeps the "Valor Mora" (F column)
# amount tied to the correct period: period 2110 keeps the value 30284
# while the rest keep 36341.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending period labels for rows 16..22
$periods = @("2104", "2105", "2106", "2107", "2108", "2109", "2110")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# Keep "Valor Mora" values attached to the right period: 2110 -> 30284,
# everything else -> 36341. Only rows 16 and 22 actually change value.
$ws.Range("F16").Value = 36341
$ws.Range("F17").Value = 36341
$ws.Range("F18").Value = 36341
$ws.Range("F19").Value = 36341
$ws.Range("F20").Value = 36341
$ws.Range("F21").Value = 36341
$ws.Range("F22").Value = 30284
